$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: "Save", formatted like the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data column H2:H4, numeric zeros, default (unstyled) formatting
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
